$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header row values: "_old" suffix -> "_FV2210", "_new" suffix -> "_FV2304"
$cols = @("A","B","C","D","E","F","G","H","I","J","L","M","N","O","P","Q","R","S","T","U")
foreach ($c in $cols) {
    $cell = $ws.Range($c + "1")
    $val = $cell.Value2
    if ($val -like "*_old") {
        $cell.Value2 = ($val -replace "_old$", "_FV2210")
    } elseif ($val -like "*_new") {
        $cell.Value2 = ($val -replace "_new$", "_FV2304")
    }
}

# 2. Freeze the header row (pane split after row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Convert range into an Excel Table (ListObject)
$rng = $ws.Range("A1:U62")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"
